$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "식당"
$ws.Range("C2").Value = 654
$ws.Range("B3").Value = "빵"
$ws.Range("C3").Value = 429
$ws.Range("B4").Value = "아침"
$ws.Range("C4").Value = 347
$ws.Range("B5").Value = "없"
$ws.Range("C5").Value = 323
$ws.Range("B6").Value = "편"
$ws.Range("C6").Value = 320
$ws.Range("B7").Value = "이용"
$ws.Range("C7").Value = 312
$ws.Range("B8").Value = "의식"
$ws.Range("C8").Value = 306
$ws.Range("B9").Value = "우유"
$ws.Range("C9").Value = 277
$ws.Range("B10").Value = "식권"
$ws.Range("C10").Value = 262
$ws.Range("B11").Value = "배식"
$ws.Range("C11").Value = 231
$ws.Range("B12").Value = "층"
$ws.Range("C12").Value = 219
$ws.Range("B13").Value = "시간"
$ws.Range("C13").Value = 214
$ws.Range("B14").Value = "식사"
$ws.Range("C14").Value = 210
$ws.Range("B15").Value = "좋"
$ws.Range("C15").Value = 207
$ws.Range("B16").Value = "사람"
$ws.Range("C16").Value = 205
$ws.Range("B17").Value = "메뉴"
$ws.Range("C17").Value = 173
$ws.Range("B18").Value = "중앙"
$ws.Range("C18").Value = 172
$ws.Range("B19").Value = "글"
$ws.Range("C19").Value = 169
$ws.Range("B20").Value = "복지관"
$ws.Range("C20").Value = 169
$ws.Range("B21").Value = "밥"
$ws.Range("C21").Value = 164
$ws.Range("B22").Value = "답변"
$ws.Range("C22").Value = 157
$ws.Range("B23").Value = "많"
$ws.Range("C23").Value = 155
$ws.Range("B24").Value = "분"
$ws.Range("C24").Value = 153
$ws.Range("B25").Value = "바나나"
$ws.Range("C25").Value = 152
$ws.Range("B26").Value = "생각"
$ws.Range("C26").Value = 147
$ws.Range("B27").Value = "식"
$ws.Range("C27").Value = 143
$ws.Range("B28").Value = "개선"
$ws.Range("C28").Value = 137
$ws.Range("B29").Value = "제공"
$ws.Range("C29").Value = 132
$ws.Range("B30").Value = "고객"
$ws.Range("C30").Value = 129
$ws.Range("B31").Value = "직원"
$ws.Range("C31").Value = 129
$ws.Range("B32").Value = "문의"
$ws.Range("C32").Value = 129
$ws.Range("B33").Value = "경우"
$ws.Range("C33").Value = 119
$ws.Range("B34").Value = "부탁"
$ws.Range("C34").Value = 110
$ws.Range("B35").Value = "감사"
$ws.Range("C35").Value = 109
$ws.Range("B36").Value = "관련"
$ws.Range("C36").Value = 109
$ws.Range("B37").Value = "안녕"
$ws.Range("C37").Value = 100
$ws.Range("B38").Value = "불편"
$ws.Range("C38").Value = 99
$ws.Range("B39").Value = "청운"
$ws.Range("C39").Value = 96
$ws.Range("B40").Value = "건의"
$ws.Range("C40").Value = 95
$ws.Range("B41").Value = "하이닉스"
$ws.Range("C41").Value = 94
$ws.Range("B42").Value = "코"
$ws.Range("C42").Value = 90
$ws.Range("B43").Value = "조식"
$ws.Range("C43").Value = 88
$ws.Range("B44").Value = "업체"
$ws.Range("C44").Value = 87
$ws.Range("B45").Value = "종류"
$ws.Range("C45").Value = 83
$ws.Range("B46").Value = "사항"
$ws.Range("C46").Value = 82
$ws.Range("B47").Value = "운영"
$ws.Range("C47").Value = 81
$ws.Range("B48").Value = "카드"
$ws.Range("C48").Value = 81
$ws.Range("B49").Value = "가능"
$ws.Range("C49").Value = 79
$ws.Range("B50").Value = "말"
$ws.Range("C50").Value = 78
$ws.Range("B51").Value = "등"
$ws.Range("C51").Value = 78
$ws.Range("B52").Value = "데"
$ws.Range("C52").Value = 77
$ws.Range("B53").Value = "사용"
$ws.Range("C53").Value = 76
$ws.Range("B54").Value = "추가"
$ws.Range("C54").Value = 76
$ws.Range("B55").Value = "줄"
$ws.Range("C55").Value = 74
$ws.Range("B56").Value = "회사"
$ws.Range("C56").Value = 74
$ws.Range("B57").Value = "식단"
$ws.Range("C57").Value = 72
$ws.Range("B58").Value = "음료"
$ws.Range("C58").Value = 70
$ws.Range("B59").Value = "주세"
$ws.Range("C59").Value = 69
$ws.Range("B60").Value = "중"
$ws.Range("C60").Value = 69
$ws.Range("B61").Value = "요청"
$ws.Range("C61").Value = 69
$ws.Range("B62").Value = "전"
$ws.Range("C62").Value = 68
$ws.Range("B63").Value = "오늘"
$ws.Range("C63").Value = 68
$ws.Range("B64").Value = "안"
$ws.Range("C64").Value = 67
$ws.Range("B65").Value = "어떻"
$ws.Range("C65").Value = 66
$ws.Range("B66").Value = "가격"
$ws.Range("C66").Value = 65
$ws.Range("B67").Value = "확인"
$ws.Range("C67").Value = 63
$ws.Range("B68").Value = "과일"
$ws.Range("C68").Value = 62
$ws.Range("B69").Value = "양"
$ws.Range("C69").Value = 61
$ws.Range("B70").Value = "샌드위치"
$ws.Range("C70").Value = 60
$ws.Range("B71").Value = "하세"
$ws.Range("C71").Value = 60
$ws.Range("B72").Value = "정도"
$ws.Range("C72").Value = 59
$ws.Range("B73").Value = "코너"
$ws.Range("C73").Value = 58
$ws.Range("B74").Value = "일"
$ws.Range("C74").Value = 56
$ws.Range("B75").Value = "선택"
$ws.Range("C75").Value = 55
$ws.Range("B76").Value = "문제"
$ws.Range("C76").Value = 55
$ws.Range("B77").Value = "이렇"
$ws.Range("C77").Value = 55
$ws.Range("B78").Value = "외부인"
$ws.Range("C78").Value = 54
$ws.Range("B79").Value = "사과"
$ws.Range("C79").Value = 54
$ws.Range("B80").Value = "면식"
$ws.Range("C80").Value = 54
$ws.Range("B81").Value = "청주"
$ws.Range("C81").Value = 53
$ws.Range("B82").Value = "안녕하"
$ws.Range("C82").Value = 53
$ws.Range("B83").Value = "외부"
$ws.Range("C83").Value = 52
$ws.Range("B84").Value = "이유"
$ws.Range("C84").Value = 52
$ws.Range("B85").Value = "앞"
$ws.Range("C85").Value = 51
$ws.Range("B86").Value = "간"
$ws.Range("C86").Value = 49
$ws.Range("B87").Value = "점심"
$ws.Range("C87").Value = 49
$ws.Range("B88").Value = "그렇"
$ws.Range("C88").Value = 49
$ws.Range("B89").Value = "세트"
$ws.Range("C89").Value = 48
$ws.Range("B90").Value = "맛"
$ws.Range("C90").Value = 48
$ws.Range("B91").Value = "이상"
$ws.Range("C91").Value = 48
$ws.Range("B92").Value = "변경"
$ws.Range("C92").Value = 47
$ws.Range("B93").Value = "요즘"
$ws.Range("C93").Value = 46
$ws.Range("B94").Value = "발생"
$ws.Range("C94").Value = 46
$ws.Range("B95").Value = "적"
$ws.Range("C95").Value = 46
$ws.Range("B96").Value = "결제"
$ws.Range("C96").Value = 46
$ws.Range("B97").Value = "공지"
$ws.Range("C97").Value = 45
$ws.Range("B98").Value = "기존"
$ws.Range("C98").Value = 45
$ws.Range("B99").Value = "안되"
$ws.Range("C99").Value = 45
$ws.Range("B100").Value = "시"
$ws.Range("C100").Value = 45
$ws.Range("B101").Value = "상태"
$ws.Range("C101").Value = 45
